# "Initial handle images #15"
#
# Populate the two previously-empty cells (C3/C4) with text, nudge the
# selection, and add the two extra picture shapes next to column C while
# re-squaring the two pre-existing pictures to their new (slightly
# tighter) bounding boxes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cell content ---------------------------------------------------
$ws.Range("C3").Value = "test"
$ws.Range("C4").Value = "test2"

# --- resize / reposition the two existing pictures -------------------
# (values are EMU/12700 = points, taken from the target bounding boxes)
$pic1 = $ws.Shapes.Item(1)
$pic1.Left = 0.028346456692913385
$pic1.Top = 0.028346456692913385
$pic1.Width = 13.266141732283465
$pic1.Height = 12.81259842519685

$pic2 = $ws.Shapes.Item(2)
$pic2.Left = 64.3748031496063
$pic2.Top = 25.568503937007875
$pic2.Width = 62.815748031496064
$pic2.Height = 24.7748031496063

# --- add the two new picture shapes near column C --------------------
$pic3 = $ws.Shapes.AddPicture("image10.png", $false, $true, 221.38582677165354, 25.2, 13.266141732283465, 12.81259842519685)
$pic4 = $ws.Shapes.AddPicture("image11.png", $false, $true, 146.5511811023622, 38.40944881889764, 13.266141732283465, 12.81259842519685)

# --- selection ---------------------------------------------------------
$ws.Range("C7").Select() | Out-Null
